# Update cryptos list price/volume columns to reflect the latest scrape
# (GitHub Actions scheduled update, Tue Aug 27 11:17:27 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cell, $text) {
    # Force the cell to stay plain text (matches the original inline-string
    # cells), preventing Excel from auto-converting numeric-looking values
    # like '154.49' into a floating point number, then drop the temporary
    # text number-format again so no stray style is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-CellText $ws.Range("D2") '62.412.88'
Set-CellText $ws.Range("D3") '2.629.06'
Set-CellText $ws.Range("E3") '  -3.88%  '
Set-CellText $ws.Range("E4") '  +0.03%  '
Set-CellText $ws.Range("D5") '551.38'
Set-CellText $ws.Range("E5") '  -2.45%  '
Set-CellText $ws.Range("D6") '154.49'
Set-CellText $ws.Range("E6") '  -4.45%  '
Set-CellText $ws.Range("E7") '  +0.11%  '
Set-CellText $ws.Range("D8") '0.587'
Set-CellText $ws.Range("E8") '  -1.59%  '
Set-CellText $ws.Range("E9") '  -4.27%  '
Set-CellText $ws.Range("E10") '  -4.17%  '
Set-CellText $ws.Range("D11") '5.44'
Set-CellText $ws.Range("E11") '  -3.40%  '
Set-CellText $ws.Range("D12") '0.363'
Set-CellText $ws.Range("E12") '  -3.32%  '
Set-CellText $ws.Range("D13") '3.094.13'
Set-CellText $ws.Range("E13") '  -3.87%  '
Set-CellText $ws.Range("D14") '25.77'
Set-CellText $ws.Range("E14") '  -4.30%  '
Set-CellText $ws.Range("D15") '62.322.70'
Set-CellText $ws.Range("E15") '  -2.08%  '
Set-CellText $ws.Range("D16") '0.0000144'
Set-CellText $ws.Range("E16") '  -3.78%  '
Set-CellText $ws.Range("D17") '2.630.69'
Set-CellText $ws.Range("E17") '  -3.95%  '
Set-CellText $ws.Range("D18") '11.68'
Set-CellText $ws.Range("E18") '  -5.13%  '
Set-CellText $ws.Range("D19") '4.54'
Set-CellText $ws.Range("E19") '  -4.16%  '
Set-CellText $ws.Range("D20") '339.90'
Set-CellText $ws.Range("E20") '  -4.41%  '
Set-CellText $ws.Range("D21") '6.09'
Set-CellText $ws.Range("E21") '  -8.28%  '
Set-CellText $ws.Range("E22") '  +0.06%  '
Set-CellText $ws.Range("D23") '0.501'
Set-CellText $ws.Range("E23") '  -3.41%  '
Set-CellText $ws.Range("D24") '62.73'
Set-CellText $ws.Range("E24") '  -2.55%  '
Set-CellText $ws.Range("E25") '  -0.88%  '
Set-CellText $ws.Range("E26") '  +0.07%  '
Set-CellText $ws.Range("D27") '8.04'
Set-CellText $ws.Range("E27") '  -3.87%  '
Set-CellText $ws.Range("D28") '0.0₃0835'
Set-CellText $ws.Range("E28") '  -8.31%  '
Set-CellText $ws.Range("E29") '  -0.13%  '
Set-CellText $ws.Range("D30") '7.03'
Set-CellText $ws.Range("E30") '  -1.82%  '
Set-CellText $ws.Range("D31") '1.89'
Set-CellText $ws.Range("E31") '  -5.09%  '
Set-CellText $ws.Range("D32") '160.56'
Set-CellText $ws.Range("E32") '  -3.76%  '
Set-CellText $ws.Range("D34") '4.75'
Set-CellText $ws.Range("E34") '  -3.35%  '
Set-CellText $ws.Range("D35") '19.22'
Set-CellText $ws.Range("E35") '  -4.18%  '
Set-CellText $ws.Range("E36") '  -4.11%  '
Set-CellText $ws.Range("D37") '1.73'
Set-CellText $ws.Range("E37") '  -4.19%  '
Set-CellText $ws.Range("D38") '334.99'
Set-CellText $ws.Range("E38") '  -3.44%  '
Set-CellText $ws.Range("D39") '6.15'
Set-CellText $ws.Range("E39") '  -2.35%  '
Set-CellText $ws.Range("D40") '0.904'
Set-CellText $ws.Range("E40") '  -7.29%  '
Set-CellText $ws.Range("D41") '3.93'
Set-CellText $ws.Range("E41") '  -3.63%  '
Set-CellText $ws.Range("D42") '37.89'
Set-CellText $ws.Range("E42") '  -2.01%  '
Set-CellText $ws.Range("D43") '0.998'
Set-CellText $ws.Range("E43") '  -0.01%  '
Set-CellText $ws.Range("D44") '20.44'
Set-CellText $ws.Range("E44") '  -6.14%  '
Set-CellText $ws.Range("D45") '0.609'
Set-CellText $ws.Range("E45") '  -3.55%  '
Set-CellText $ws.Range("E46") '  -0.67%  '
Set-CellText $ws.Range("D47") '19.68'
Set-CellText $ws.Range("E47") '  -6.73%  '
Set-CellText $ws.Range("D48") '0.0547'
Set-CellText $ws.Range("E48") '  -6.15%  '
Set-CellText $ws.Range("D49") '0.0961'
Set-CellText $ws.Range("E49") '  -3.49%  '
Set-CellText $ws.Range("D50") '127.60'
Set-CellText $ws.Range("E50") '  -3.77%  '
Set-CellText $ws.Range("D51") '0.0238'
Set-CellText $ws.Range("E51") '  -4.71%  '
